$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add a new forecast column (V) + a new "as-of 2020-05-15" row (34) to both
# the "cases" and "deaths" sheets, following the same staircase-forecast
# layout already used by the table. The new row label (A34) carries the new
# "2020-05-15" date text; the new column header (V1), per source data,
# re-uses the existing "2020-05-01" label already used elsewhere in the
# sheet (matching the exact target content, not a "corrected" label).
# ---------------------------------------------------------------------------

# Prefixing with an apostrophe forces Excel to store these as literal text
# (matching the existing "yyyy-mm-dd" text labels) instead of auto-
# converting them to date serial numbers.
$newColHeaderLabel = "'2020-05-01"
$newRowLabel = "'2020-05-15"

# New column is V (22nd column); new row is 34.
$newCol = 22
$newRow = 34

# Per-sheet numbers: B20 (previously blank) gets a value, and the new
# column V gets the staircase of forecast values for rows 21-34.
$sheetData = @{
    "cases"  = @{
        B20 = 30374
        V   = @{
            21 = 32739
            22 = 35440
            23 = 38033
            24 = 41669
            25 = 45105
            26 = 48495
            27 = 51479
            28 = 54087
            29 = 55793
            30 = 58460
            31 = 60793
            32 = 62978
            33 = 64822
            34 = 66346
        }
    }
    "deaths" = @{
        B20 = 2511
        V   = @{
            21 = 2672
            22 = 2895
            23 = 2997
            24 = 3284
            25 = 3455
            26 = 3611
            27 = 3757
            28 = 3873
            29 = 3947
            30 = 4078
            31 = 4181
            32 = 4253
            33 = 4344
            34 = 4403
        }
    }
}

foreach ($sheetName in @("cases", "deaths")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $sheetData[$sheetName]

    # --- Header: V1 = "2020-05-01" label (re-used, per source data) -----
    $cell = $ws.Cells.Item(1, $newCol)
    $cell.Value = $newColHeaderLabel
    $cell.Style = "Normal"

    # --- Existing rows 2-19: add a blank V cell to extend the table -----
    for ($r = 2; $r -le 19; $r++) {
        $c = $ws.Cells.Item($r, $newCol)
        $c.NumberFormat = "General"
        $c.Style = "Normal"
    }

    # --- Row 20: B20 gains a value, V20 stays blank (new, empty) --------
    $ws.Cells.Item(20, 2).Value = $info.B20

    $c = $ws.Cells.Item(20, $newCol)
    $c.NumberFormat = "General"
    $c.Style = "Normal"

    # --- Rows 21-33: add the new diagonal forecast value in column V ----
    for ($r = 21; $r -le 33; $r++) {
        $ws.Cells.Item($r, $newCol).Value = $info.V[$r]
    }

    # --- New row 34: label in A34, blanks in B34:U34, value in V34 ------
    $c = $ws.Cells.Item($newRow, 1)
    $c.Value = $newRowLabel
    $c.Style = "Normal"

    for ($col = 2; $col -le 21; $col++) {
        $c = $ws.Cells.Item($newRow, $col)
        $c.NumberFormat = "General"
        $c.Style = "Normal"
    }

    $ws.Cells.Item($newRow, $newCol).Value = $info.V[$newRow]
}
